$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.025.80'
$ws.Range('E2').Value = '  +4.49%  '
$ws.Range('D3').Value = '2.742.45'
$ws.Range('E3').Value = '  +3.42%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.74'
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '159.61'
$ws.Range('E6').Value = '  +10.67%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.621'
$ws.Range('E7').Value = '  +3.72%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').Value = '2.760.22'
$ws.Range('E9').Value = '  +3.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.85'
$ws.Range('E10').Value = '  +4.06%  '
$ws.Range('E11').Value = '  +2.40%  '
$ws.Range('E12').Value = '  +3.19%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.160'
$ws.Range('E13').Value = '  +0.85%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.83'
$ws.Range('E14').Value = '  +6.53%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '3.226.58'
$ws.Range('E15').Value = '  +3.35%  '
$ws.Range('D16').Value = '63.710.63'
$ws.Range('E16').Value = '  +4.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000157'
$ws.Range('E17').Value = '  +7.16%  '
$ws.Range('D18').Value = '2.752.35'
$ws.Range('E18').Value = '  +3.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.18'
$ws.Range('E19').Value = '  +4.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.97'
$ws.Range('E20').Value = '  +3.91%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '363.67'
$ws.Range('E21').Value = '  +2.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.02'
$ws.Range('E22').Value = '  +2.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.544'
$ws.Range('E23').Value = '  +3.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.995'
$ws.Range('E24').Value = '  -0.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '67.11'
$ws.Range('E25').Value = '  +3.86%  '
$ws.Range('E26').Value = '  +5.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.61'
$ws.Range('E27').Value = '  +1.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('D29').Value = '0.0₃0926'
$ws.Range('E29').Value = '  +13.20%  '
$ws.Range('E30').Value = '  +1.22%  '
$ws.Range('E31').Value = '  +5.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.27'
$ws.Range('E32').Value = '  +12.69%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '173.88'
$ws.Range('E33').Value = '  +2.90%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '20.66'
$ws.Range('E35').Value = '  +2.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.97'
$ws.Range('E36').Value = '  +6.27%  '
$ws.Range('E37').Value = '  +5.94%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.83'
$ws.Range('E38').Value = '  +6.09%  '
$ws.Range('E39').Value = '  +3.98%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.26'
$ws.Range('E40').Value = '  +2.43%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.25'
$ws.Range('E41').Value = '  +16.33%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '339.26'
$ws.Range('E42').Value = '  +0.37%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '39.61'
$ws.Range('E43').Value = '  +2.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.42'
$ws.Range('E44').Value = '  +6.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.85'
$ws.Range('E45').Value = '  +5.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0602'
$ws.Range('E46').Value = '  +3.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.644'
$ws.Range('E47').Value = '  +2.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0260'
$ws.Range('E48').Value = '  +2.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '137.23'
$ws.Range('E49').Value = '  +1.19%  '
$ws.Range('E50').Value = '  +2.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.996'
$ws.Range('E51').Value = '  -0.20%  '
